$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (3:22) down by one row (to 4:23), working bottom-up
# so we never clobber a row before it has been read.
for ($r = 22; $r -ge 3; $r--) {
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r + 1, $c).Value = $ws.Cells.Item($r, $c).Value2
    }
}

# Populate the newly freed row 3 with the new observation (2020-04-01).
$ws.Cells.Item(3, 1).Value = "2020-04-01 00:00:00_diff"
$ws.Cells.Item(3, 2).Value = 9.006385916832064
$ws.Cells.Item(3, 3).Value = -8.778729625081095
$ws.Cells.Item(3, 4).Value = -0.3907511347594377
$ws.Cells.Item(3, 5).Value = 1.872451415860664
$ws.Cells.Item(3, 6).Value = -1.508723736095618
$ws.Cells.Item(3, 7).Value = -1.540727864400714
$ws.Cells.Item(3, 8).Value = 0.5152984556711749

# Row 23 is a brand-new row (it didn't exist before), so it has no inherited
# formatting yet -- copy the label style (bold / bordered / centered) used by
# the rest of column A onto it.
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(23, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
